# Add two new header columns ("gender" in E1, "dob" in F1) to the
# student upload template, matching the header style of the existing
# columns (A1:D1), and update the active selection to G8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers. Write F1 ("dob") before E1 ("gender") so the shared
# strings table picks up "dob" at index 4 and "gender" at index 5,
# matching the target workbook.
$ws.Range("F1").Value = "dob"
$ws.Range("E1").Value = "gender"

# Match the bold header formatting used by the existing header cells.
$ws.Range("E1").Font.Bold = $true
$ws.Range("F1").Font.Bold = $true

# Update the selected cell/active cell as recorded in the saved view.
[void]$ws.Range("G8").Select()
